$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "241.85"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "21.54"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.273"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05606"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.383"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.376"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8073"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9441"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01115"
$ws.Range("E10").Value = "9OneONE"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1425"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07418"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03235"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03064"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09277"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.586"
$ws.Range("E16").Value = "15MCDexMCB"
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.001657"
$ws.Range("E17").Value = "16BitForexTokenBF"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.04710"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006367"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004971"
$ws.Range("B21").Value = "UpBots"
$ws.Range("C21").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.007505"
$ws.Range("E21").Value = "20UpBotsUBXTBestin24h"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.001044"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0001503"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.769"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.103"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.3253"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03901"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006734"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1030"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003114"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007528"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005958"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000752"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005506"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.6841"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05625"
$ws.Range("E49").Value = "48BOLOBOLO"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002105"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.01012"
